# server commands.xlsx - add checkpoint command, clarify dump arg description
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clarify the existing "dump" command's argument description.
$ws.Range("E7").Value = "target board address (int)"

# New row: checkpoint command
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "checkpoint"
$ws.Range("C8").Value = "creates savefile checkpoint"
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = "runname (str)"

# Update selection to match the author's final cursor position.
$ws.Range("D10").Select()
